# Updated cryptos list on Fri Mar 24 07:57:18 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.327.67"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.30%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.817.66"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.60%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "325.47"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9991"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4346"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -3.60%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3679"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.82%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "44.89"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.77%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07677"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.28%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.148"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.19%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "22.03"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.49%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.328"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.44%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.494"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.80%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.817.77"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.76%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "95.37"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +8.37%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001080"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06444"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.62%  "
$ws.Range("E20").Value = "  -0.06%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.45"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.80%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.246"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.22%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "28.337.29"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.26%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.60"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.142"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -7.69%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "159.82"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.22%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.73"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.60%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.026.78"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.94%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.281"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -3.31%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "131.79"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +3.03%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.201"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.43%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.022"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.03%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.09136"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.86%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.560"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.14%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "13.03"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.98%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02412"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +4.28%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.226"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.53%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2177"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.16%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.6604"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.92%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.06205"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.46%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.203"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("E42").Value = "  +1.04%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.430"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.95%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.9987"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.89"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.71%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.6111"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.77%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.736"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.57%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "126.01"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.022"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.20%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.166"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.46%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06997"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.45%  "
